$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 182
$wsOff.Range("C2").Value = 131
$wsOff.Range("D2").Value = 40
$wsOff.Range("E2").Value = 16

# Sheet "DEF" - row 2 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 189
$wsDef.Range("C2").Value = 124
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 19
$wsDef.Range("F2").Value = 2
